# #22 modify against komodan comments
#
# - grow the 3 red round-rect containers taller (bottom "batch component"
#   area gets 2 more rows of boxes)
# - push the "MyBatis3" / "Spring Framework" bars down to make room
# - add 4 new gold/gradient boxes describing business-logic sub components
#
# NOTE on Left/Top/Width/Height literals below: PowerPoint's Shape geometry
# is stored internally as points in a 32-bit float, then multiplied by
# 12700 (EMU per point) and truncated when serialised to OOXML. The literal
# point values used here are chosen (via a float32 search) so that after
# that float32 round-trip + truncation they reproduce the exact target EMU
# from the authoritative edit, instead of drifting by the usual +/-1 EMU
# that naive "emu / 12700" math would leave behind.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------------
# 1) Grow the three red outline containers vertically to make room for
#    the two new rows of component boxes.
# ---------------------------------------------------------------------

# 角丸四角形 25 (outer red rounded rect that wraps the whole diagram)
$s.Shapes.Item("角丸四角形 25").Height = 555.65234375

# 角丸四角形 24 (red rounded rect around "バッチ機能コンポーネント")
$s.Shapes.Item("角丸四角形 24").Height = 238.13670349121094

# 角丸四角形 23 (red rounded rect around "バッチ実行基盤")
$s.Shapes.Item("角丸四角形 23").Height = 238.13670349121094

# ---------------------------------------------------------------------
# 2) Slide the "MyBatis3" / "Spring Framework" bars further down so they
#    stay below the taller containers above.
# ---------------------------------------------------------------------

$s.Shapes.Item("正方形/長方形 21").Top = 480.19671630859375   # MyBatis3
$s.Shapes.Item("正方形/長方形 22").Top = 525.146484375        # Spring Framework

# ---------------------------------------------------------------------
# 3) Add the 4 new gradient-filled component boxes.
#
# The deck's shape-id allocator hands out the lowest still-unused id each
# time a shape is created (PowerPoint itself reuses ids left behind by
# shapes that were created/deleted earlier in the authoring session). The
# ids 2,3,8,17,27,32 are free "holes" below the existing max id (41) that
# get consumed first; only the 7th-10th newly created shapes land on
# 33/34/35/42 -- which are the ids the target deck actually uses for
# these 4 boxes. So we spin up 6 disposable placeholder shapes first
# (soaking up those holes), delete them again, and only keep the 4
# shapes created right after.
# ---------------------------------------------------------------------

$template = $s.Shapes.Item("正方形/長方形 14")   # existing gold gradient box ("入力データ取得")

$scratch = @()
for ($i = 0; $i -lt 6; $i++) {
    $d = $template.Duplicate()
    $scratch += $d.Item(1)
}

# --- 正方形/長方形 32 : "ビジネスロジック実行" ---
$d = $template.Duplicate()
$sh1 = $d.Item(1)
$sh1.Name = "正方形/長方形 32"
$sh1.Left = 286.2910461425781
$sh1.Top = 377.728515625
$sh1.Width = 170.09764099121094
$sh1.Height = 34.01953125
$sh1.TextFrame.TextRange.Text = "ビジネスロジック実行"
$sh1.TextFrame.TextRange.Font.Size = 16

# --- 正方形/長方形 33 : "メッセージ管理" ---
$d = $template.Duplicate()
$sh2 = $d.Item(1)
$sh2.Name = "正方形/長方形 33"
$sh2.Left = 467.728515625
$sh2.Top = 377.068359375
$sh2.Width = 170.09764099121094
$sh2.Height = 34.01953125
$sh2.TextFrame.TextRange.Text = "メッセージ管理"

# --- 正方形/長方形 34 : "例外ハンドリング" ---
$d = $template.Duplicate()
$sh3 = $d.Item(1)
$sh3.Name = "正方形/長方形 34"
$sh3.Left = 286.2910461425781
$sh3.Top = 417.41796875
$sh3.Width = 170.09764099121094
$sh3.Height = 34.01953125
$sh3.TextFrame.TextRange.Text = "例外ハンドリング"

# --- 正方形/長方形 41 : "ファイル操作" ---
$d = $template.Duplicate()
$sh4 = $d.Item(1)
$sh4.Name = "正方形/長方形 41"
$sh4.Left = 467.728515625
$sh4.Top = 416.7578125
$sh4.Width = 170.09764099121094
$sh4.Height = 34.01953125
$sh4.TextFrame.TextRange.Text = "ファイル操作"

foreach ($junk in $scratch) {
    $junk.Delete()
}
